# ND03.xlsx update: "no blank columns for now"
#
# 1. WMT_Extract (sheet 1): the report-code label cells change from "ND03" to "ND01".
# 2. Inst_Reports (sheet 3): a previously-blank G2 cell is filled in with 2 (so row 2
#    matches row 3's F/G pairing instead of leaving a gap between F2 and H2).
# 3. View/selection bookkeeping: Inst_Reports becomes the active/selected tab instead
#    of WMT_Extract, and each sheet's remembered selection moves on one cell.

$wb = $excel.ActiveWorkbook

$wsExtract = $wb.Worksheets.Item("WMT_Extract")
$wsInst    = $wb.Worksheets.Item("Inst_Reports")

# -- Data edits --------------------------------------------------------------

# Report code was "ND03", should read "ND01".
$wsExtract.Range("C2").Value = "ND01"
$wsExtract.Range("C3").Value = "ND01"

# Fill the gap at G2 (was skipped between F2 and H2) with 2, same as G3.
$wsInst.Range("G2").Value = 2

# -- View / selection edits --------------------------------------------------

# WMT_Extract is no longer the selected tab; its remembered selection moves
# from AE1 to AE2.
$wsExtract.Range("AE2").Select()

# Inst_Reports becomes the active (selected) tab, with its selection moved
# from D26 to E14.
$wsInst.Activate()
$wsInst.Range("E14").Select()
